# Updating todo list for panzer infrastructure
#
# 1. Fix "phsyics" -> "physics" typo (two occurrences in one paragraph).
# 2. Mark the "Split the registerGatherScatter into:" bullet as [FINISHED]
#    (in red) while leaving the rest of the sentence alone.
# 3. Append a long red "UPDATE (2012.05.25): ..." note to the
#    "Eliminate BCStrategy..." bullet.
# 4. Turn the numbered-list levels 3-9 (ilvl 2-8) into fully qualified
#    multi-level numbers (e.g. %2.%3. instead of just %3.).
# 5. Recolor the Normal style's default font color.
# 6. Merge "Numbering Symbols" / "Bullets" character styles into a single
#    "Bullets" style and add two new "ListLabel" character styles.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Typo fix: "phsyics" -> "physics" in the volume/side physics block
#    paragraph only (there is another, unrelated, "Unique phsyics block"
#    paragraph elsewhere that must stay untouched).
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*A volume phsyics block is cloned*") {
        $p.Range.Find.Execute("phsyics", $false, $false, $false, $false, $false, $true, 1, $false, "physics", 2)
    }
}

# ---------------------------------------------------------------------
# 2) "[FINISHED] " prefix, in red, before "Split the registerGatherScatter..."
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Split the registerGatherScatter into:*") {
        $insPoint = $p.Range.Duplicate
        $insPoint.Collapse(1)
        $insPoint.InsertBefore("[FINISHED] ")

        $redRange = $p.Range.Duplicate
        $redRange.SetRange($p.Range.Start, $p.Range.Start + 10)
        $redRange.Font.Color = 255
        break
    }
}

# ---------------------------------------------------------------------
# 3) Append ".  " then a new red "UPDATE (2012.05.25): ..." run to the
#    "Eliminate BCStrategy..." bullet (same paragraph, two runs).
# ---------------------------------------------------------------------
$updateText = "UPDATE (2012.05.25):  I don't think this is possible.  Since bcs span physics from multiple equation sets, it needs the physics block to build certain evaluators.  However the interface for equation sets knows nothing about the physics block since the physics block own multiple equation sets.  We would have to modify the equation sets to accept the physics block as an argument which is a little strange and doesn't make sense for volume assembly.  We will have to think more about this.  Most of the machinery for both equation sets and bc_strategies is the same so maybe we build some basic tools that both can specialize from.  But the more I think about it, I think the generalized part of the tools already exist.  We may have to live with this for now.  It might take a really big redesign to address this."

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Eliminate BCStrategy*") {
        $r = $p.Range
        $tail = $d.Range($r.End - 1, $r.End - 1)
        $tail.InsertAfter(".  ")

        $r2 = $p.Range
        $ins2 = $d.Range($r2.End - 1, $r2.End - 1)
        $updateStart = $ins2.Start
        $ins2.InsertAfter($updateText)

        $redRange = $d.Range($updateStart, $updateStart + $updateText.Length)
        $redRange.Font.Color = 255
        break
    }
}

# ---------------------------------------------------------------------
# 4) Multi-level numbering: ilvl 2-8 (levels 3-9) should concatenate the
#    ancestor level numbers, e.g. "%3." -> "%2.%3.", "%9." -> "%2.%3.%4.%5.%6.%7.%8.%9."
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Split the registerGatherScatter into:*") {
        $lt = $p.Range.ListFormat.ListTemplate
        $lt.ListLevels(3).NumberFormat = "%2.%3."
        $lt.ListLevels(4).NumberFormat = "%2.%3.%4."
        $lt.ListLevels(5).NumberFormat = "%2.%3.%4.%5."
        $lt.ListLevels(6).NumberFormat = "%2.%3.%4.%5.%6."
        $lt.ListLevels(7).NumberFormat = "%2.%3.%4.%5.%6.%7."
        $lt.ListLevels(8).NumberFormat = "%2.%3.%4.%5.%6.%7.%8."
        $lt.ListLevels(9).NumberFormat = "%2.%3.%4.%5.%6.%7.%8.%9."
        break
    }
}

# ---------------------------------------------------------------------
# 5) Normal style's default run color: auto -> 00000A
#    (Word's Color longs are 0x00BBGGRR, so RGB 00000A -> 0x0A0000)
# ---------------------------------------------------------------------
$normal = $d.Styles("Normal")
$normal.Font.Color = 655360

# ---------------------------------------------------------------------
# 6) Style table surgery:
#    - drop the duplicate "Bullets" character style
#    - rename "Numbering Symbols" to "Bullets" and give it the OpenSymbol
#      rFonts that used to live on the (now deleted) duplicate
#    - add two new "ListLabel 1" / "ListLabel 2" character styles
# ---------------------------------------------------------------------
$oldBullets = $d.Styles("Bullets")
$oldBullets.Delete()

$numSym = $d.Styles("Numbering Symbols")
$numSym.NameLocal = "Bullets"
$numSym.Font.NameAscii = "OpenSymbol"
$numSym.Font.NameFarEast = "OpenSymbol"
$numSym.Font.NameOther = "OpenSymbol"
$numSym.Font.NameBi = "OpenSymbol"

$listLabel1 = $d.Styles.Add("ListLabel 1", 2)
$listLabel1.Font.NameBi = "Wingdings 2"
$listLabel1.NextParagraphStyle = "ListLabel 1"

$listLabel2 = $d.Styles.Add("ListLabel 2", 2)
$listLabel2.Font.NameBi = "OpenSymbol"
$listLabel2.NextParagraphStyle = "ListLabel 2"
